$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B6").Value = "15-18"
$ws.Range("B9").Select()
